# Insert a new data row at row 428 (shifting all subsequent rows down by one,
# which also matches the rest of the diff where every row from the old 428..526
# reappears, unchanged, as the new 429..527).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(428).Insert()

# Populate the newly inserted row 428 with its data.
$ws.Range("A428").Value = 10
$ws.Range("B428").Value = "Vega Modelo de Temuco"
$ws.Range("C428").Value = "La Araucanía"
$ws.Range("D428").Value = (Get-Date -Year 2023 -Month 9 -Day 4 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("D428").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E428").Value = 9
$ws.Range("F428").Value = 100114013
$ws.Range("G428").Value = "Zanahoria"
$ws.Range("H428").Value = "Sin especificar"
$ws.Range("I428").Value = "Primera"
$ws.Range("J428").Value = 120
$ws.Range("K428").Value = 5000
$ws.Range("L428").Value = 5000
$ws.Range("M428").Value = 5000
$ws.Range("N428").Value = "`$/saco 20 kilos"
$ws.Range("O428").Value = "Región de La Araucanía"
$ws.Range("P428").Value = 250
$ws.Range("Q428").Value = 20
$ws.Range("R428").Value = "Hortaliza"
